# 200212 animation add, map lock add, character lock add, prefab modified
#
# B8 flips from the "최준아" placeholder to a completed "O" mark, and a brand
# new "O" mark is added at B22 (the row that previously had no status yet).
# B6/B7 keep showing "최준아" (no visible change there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "O"
$ws.Range("B22").Value = "O"

# Move the on-screen selection/scroll position the way the author last left it.
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
